$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Years_of_data formulas to add 1 (inclusive year count)
$ws.Range("Q2").Formula = "=N2-K2+1"
$ws.Range("Q3").Formula = "=N3-K3+1"
$ws.Range("Q4").Formula = "=N4-K4+1"

# Update the active cell selection
$ws.Range("Q5").Select()
